$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values differ between row 4 and row 5 and must be swapped.
$cols = @(1, 2, 4, 5, 6, 7, 8, 9, 17, 18)   # A, B, D, E, F, G, H, I, Q, R

foreach ($col in $cols) {
    $cell4 = $ws.Cells.Item(4, $col)
    $cell5 = $ws.Cells.Item(5, $col)
    $v4 = $cell4.Value2
    $v5 = $cell5.Value2
    $cell4.Value2 = $v5
    $cell5.Value2 = $v4
}

# AO4 ("brandstubbe") moves to AO5; AO4 becomes empty.
$ao4 = $ws.Cells.Item(4, 41)
$ao5 = $ws.Cells.Item(5, 41)
$ao5.Value2 = $ao4.Value2
$ao4.Value2 = ""
